$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day 7 follow-up dictionary: the "hf_id" label (column B, row 7) is
# renamed to "fid".
$ws.Range("B7").Value = "fid"

# Restore the saved view state: scroll so row 2 is the top visible row and
# move the active selection to D15.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("D15").Select() | Out-Null
